$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D31").Value = 284.93
$ws1.Range("L31").Value = 1108.6
$ws1.Range("D60").Value = "5 de 58"
$ws1.Range("L60").Value = "7 de 58"

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F31").Value = 1400.66
$ws2.Range("F60").Value = 13223.01

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 4210.36
$ws3.Range("E3").Value = 10615.05
$ws3.Range("F3").Value = 0.283996193022655

$ws3.Range("D11").Value = 3838.25
$ws3.Range("E11").Value = 12309.75
$ws3.Range("F11").Value = 0.2376919742382958

$ws3.Range("D12").Value = 5645.73
$ws3.Range("E12").Value = 44661.27
$ws3.Range("F12").Value = 0.1122255352137873

$ws3.Range("D14").Value = 14795.14
$ws3.Range("E14").Value = 83066.74766749098
$ws3.Range("F14").Value = 0.1511838812088931
